$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B-E in this sheet store plain text (prices, URLs, percentages),
# even when the text looks like a pure number (e.g. "4.93" or "0.370").
# Force text format on column D before assignment so Excel does not
# auto-convert these into numeric values (which would also drop trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "75.736.68"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "2.921.23"
$ws.Range("E3").Value = "  +4.79%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "200.99"
$ws.Range("E5").Value = "  +7.54%  "
$ws.Range("D6").Value = "596.11"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.553"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").Value = "0.196"
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("D10").Value = "2.926.70"
$ws.Range("E10").Value = "  +5.09%  "
$ws.Range("D11").Value = "0.445"
$ws.Range("E11").Value = "  +19.42%  "
$ws.Range("D12").Value = "0.161"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "4.93"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").Value = "3.466.40"
$ws.Range("E14").Value = "  +4.71%  "
$ws.Range("D15").Value = "28.18"
$ws.Range("E15").Value = "  +5.80%  "
$ws.Range("D16").Value = "75.782.60"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "0.0000189"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("D18").Value = "2.929.01"
$ws.Range("E18").Value = "  +4.80%  "
$ws.Range("D19").Value = "13.19"
$ws.Range("E19").Value = "  +8.55%  "
$ws.Range("D20").Value = "8.75"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").Value = "373.09"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "2.30"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").Value = "4.32"
$ws.Range("E23").Value = "  +6.40%  "
$ws.Range("D24").Value = "71.93"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").Value = "3.090.56"
$ws.Range("E26").Value = "  +4.89%  "
$ws.Range("D27").Value = "4.31"
$ws.Range("E27").Value = "  +5.09%  "
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "0.0000108"
$ws.Range("E29").Value = "  +6.33%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "1.38"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "7.88"
$ws.Range("E32").Value = "  +4.55%  "
$ws.Range("D33").Value = "499.83"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E34").Value = "  +3.94%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "20.23"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("D37").Value = "163.78"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").Value = "0.109"
$ws.Range("E38").Value = "  +26.96%  "
$ws.Range("D39").Value = "19.63"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").Value = "0.370"
$ws.Range("E40").Value = "  +9.32%  "
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  -3.55%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "178.79"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "4.97"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "1.65"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").Value = "40.18"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "1.19"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "0.577"
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "3.84"
$ws.Range("E50").Value = "  +4.27%  "
$ws.Range("D51").Value = "22.53"
$ws.Range("E51").Value = "  +9.10%  "
